# Applies the cryptos.xlsx price/volume update described in the commit
# "Updated cryptos list on Wed Sep  6 18:58:13 UTC 2023 with GitHub Actions".
#
# Column D ("Price") cells are plain text in this workbook (prices use "."
# as a thousands separator in some rows, e.g. "25.757.13"), so writing a
# value that looks like a simple decimal (e.g. "215.49") through .Value
# would otherwise be auto-coerced to a numeric cell by Excel. To keep those
# cells as text we briefly force a text NumberFormat before assigning the
# value, then restore the cell style so formatting is left untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "25.757.13"
$ws.Range("E2").Value = "  -0.20%  "

Set-TextValue $ws.Range("D3") "1.636.39"
$ws.Range("E3").Value = "  -0.20%  "

$ws.Range("E4").Value = "  -0.06%  "

Set-TextValue $ws.Range("D5") "215.49"
$ws.Range("E5").Value = "  -0.08%  "

$ws.Range("E6").Value = "  -0.79%  "

$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("E8").Value = "  -0.09%  "

$ws.Range("E9").Value = "  -0.95%  "

Set-TextValue $ws.Range("D10") "19.58"
$ws.Range("E10").Value = "  -3.97%  "

$ws.Range("E11").Value = "  +1.00%  "

$ws.Range("E12").Value = "  -0.51%  "

$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-TextValue $ws.Range("D13") "1.860.96"
$ws.Range("E13").Value = "  -0.18%  "

$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue $ws.Range("D14") "1.633.14"
$ws.Range("E14").Value = "  -0.32%  "

Set-TextValue $ws.Range("D15") "0.556"
$ws.Range("E15").Value = "  -0.97%  "

Set-TextValue $ws.Range("D16") "0.0₃0768"
$ws.Range("E16").Value = "  +0.35%  "

Set-TextValue $ws.Range("D17") "62.87"
$ws.Range("E17").Value = "  -0.78%  "

Set-TextValue $ws.Range("D18") "25.779.06"
$ws.Range("E18").Value = "  -0.15%  "

$ws.Range("E19").Value = "  +0.02%  "

$ws.Range("E20").Value = "  +1.39%  "

Set-TextValue $ws.Range("D21") "193.69"
$ws.Range("E21").Value = "  +0.48%  "

Set-TextValue $ws.Range("D22") "9.97"
$ws.Range("E22").Value = "  +0.36%  "

$ws.Range("E23").Value = "  +2.15%  "

$ws.Range("E24").Value = "  +0.04%  "

$ws.Range("E25").Value = "  +2.27%  "

Set-TextValue $ws.Range("D26") "140.31"
$ws.Range("E26").Value = "  -0.49%  "

$ws.Range("E27").Value = "  -1.97%  "

$ws.Range("E28").Value = "  +0.74%  "

Set-TextValue $ws.Range("D29") "15.52"
$ws.Range("E29").Value = "  -0.45%  "

$ws.Range("E30").Value = "  -0.12%  "

$ws.Range("E31").Value = "  -0.12%  "

$ws.Range("E32").Value = "  +1.63%  "

$ws.Range("E34").Value = "  +1.34%  "

$ws.Range("E35").Value = "  +0.24%  "

Set-TextValue $ws.Range("D36") "0.898"
$ws.Range("E36").Value = "  -0.57%  "

$ws.Range("E37").Value = "  -1.12%  "

$ws.Range("B38").Value = "MXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue $ws.Range("D38") "2.52"
$ws.Range("E38").Value = "  -1.93%  "

$ws.Range("B39").Value = "Maker"
$ws.Range("C39").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue $ws.Range("D39") "1.117.06"
$ws.Range("E39").Value = "  -1.39%  "

$ws.Range("E40").Value = "  -0.38%  "

$ws.Range("E41").Value = "  +0.65%  "

$ws.Range("E42").Value = "  +1.47%  "

Set-TextValue $ws.Range("D43") "99.68"
$ws.Range("E43").Value = "  +0.74%  "

$ws.Range("E44").Value = "  -0.07%  "

Set-TextValue $ws.Range("D45") "1.770.66"

Set-TextValue $ws.Range("D46") "0.0₆0110"
$ws.Range("E46").Value = "  +1.18%  "

Set-TextValue $ws.Range("D47") "55.23"
$ws.Range("E47").Value = "  -0.68%  "

$ws.Range("E48").Value = "  -2.41%  "

$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D49") "7.64"
$ws.Range("E49").Value = "  -1.71%  "

$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws.Range("D50") "0.0501"
$ws.Range("E50").Value = "  -0.40%  "

$ws.Range("E51").Value = "  +2.81%  "
